# Updated legacy GSC export data.
#
# The Coverage ("Chart") sheet's first data row (2025-11-16) is dropped, so
# every remaining date (2025-11-17 .. 2026-02-10) keeps the Not indexed /
# Indexed / Impressions figures that used to belong to the following day,
# and the final day (which no longer has a "next day" to inherit from)
# drops off the bottom of the table.
#
# Deleting row 2 (the 2025-11-16 row) and letting Excel shift everything
# else up reproduces exactly that transformation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$ws.Rows.Item(2).Delete()
